$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Label" header in H1, matching the header style used by G1 (bold, centered, bordered)
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

# --- Block for 100 iterations (rows 2-7): refit D/E/F values + new Label (H) column ---
$ws.Range("D2").Value = 0.6217755915230598
$ws.Range("E2").Value = 0.6217755915230598
$ws.Range("H2").Value = 0

$ws.Range("D3").Value = 0.6769681701771864
$ws.Range("E3").Value = 0.6769681701771864
$ws.Range("H3").Value = 0

$ws.Range("D4").Value = 0.4912260747114239
$ws.Range("E4").Value = 0.4912260747114239
$ws.Range("H4").Value = 0

$ws.Range("H5").Value = 0

$ws.Range("D6").Value = 0.4731063092683897
$ws.Range("E6").Value = 0.4731063092683897
$ws.Range("H6").Value = 0

$ws.Range("D7").Value = 0.5305438823085764
$ws.Range("E7").Value = 0.4694561176914236
$ws.Range("F7").Value = 0.9024326801300049
$ws.Range("H7").Value = 1

# --- Block for 200 iterations (rows 8-13): only the new Label (H) column changes ---
$ws.Range("H8").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("H10").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 1
